# Updated cryptos list on Thu Nov 23 13:32:26 UTC 2023 with GitHub Actions
# Refreshes price/volume(1h) figures for each coin row, and reflects the
# new rank ordering for the ImmutableX/Cosmos and TrustWalletToken/Cronos pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.362.50"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "2.060.29"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'233.44"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +3.17%  "
$ws.Range("D7").Value = "'57.83"
$ws.Range("E7").Value = "  +5.42%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +3.70%  "
$ws.Range("D10").Value = "'58.31"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D13").Value = "2.367.37"
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "'14.40"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "'20.76"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").Value = "'5.14"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "2.054.18"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "37.322.17"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").Value = "'6.36"
$ws.Range("E20").Value = "  +18.00%  "
$ws.Range("D21").Value = "'69.25"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("D23").Value = "'225.27"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").Value = "'2.40"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").Value = "'166.28"
$ws.Range("E27").Value = "  +2.62%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'1.46"
$ws.Range("E28").Value = "  +6.55%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'8.82"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").Value = "'0.129"
$ws.Range("E30").Value = "  -3.07%  "
$ws.Range("D31").Value = "'19.10"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").Value = "'0.118"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").Value = "'4.49"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").Value = "'2.56"
$ws.Range("E34").Value = "  +4.21%  "
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("D36").Value = "'4.54"
$ws.Range("E36").Value = "  +7.16%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'5.89"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "'3.27"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "'4.65"
$ws.Range("E41").Value = "  +12.47%  "
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").Value = "1.477.69"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("D44").Value = "'96.86"
$ws.Range("E44").Value = "  +7.27%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  +6.12%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.0926"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("D48").Value = "'15.55"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("D50").Value = "'7.17"
$ws.Range("E50").Value = "  +4.52%  "
$ws.Range("D51").Value = "'2.95"
$ws.Range("E51").Value = "  +2.15%  "
